# Update demo_regions as well
# Applies the data-entry updates made across several "Programs ..." input
# sheets, a row-height tweak on "Relative risks", and moves the active
# sheet/tab to "Programs target population".

$wb = $excel.ActiveWorkbook

# --- Relative risks: row 56 is now a taller (wrapped) row ---
$wsRelRisks = $wb.Worksheets.Item("Relative risks")
$wsRelRisks.Rows.Item(56).RowHeight = 26

# --- Programs birth outcomes: IFAS for pregnant women effectiveness -> 0 ---
$wsBirth = $wb.Worksheets.Item("Programs birth outcomes")
$wsBirth.Activate()
$wsBirth.Range("C4:D4").Value = 0
$wsBirth.Range("C6:D6").Value = 0
$wsBirth.Range("C8").Select()

# --- Programs anemia: coverage figures updated to 0.976 ---
$wsAnemia = $wb.Worksheets.Item("Programs anemia")
$wsAnemia.Activate()
$wsAnemia.Range("E20:O20").Value = 0.976
$wsAnemia.Range("E19:O20").Select()

# --- Programs wasting: odds ratios updated ---
$wsWasting = $wb.Worksheets.Item("Programs wasting")
$wsWasting.Activate()
$wsWasting.Range("D3:G3").Value = 0.22
$wsWasting.Range("D5:G5").Value = 0.16
$wsWasting.Range("D5:G5").Select()

# --- Programs for children: effectiveness / affected-fraction values updated ---
$wsChildren = $wb.Worksheets.Item("Programs for children")
$wsChildren.Activate()
$wsChildren.Range("F3:H3").Value = 0.36
$wsChildren.Range("F4:H4").Value = 0.45
$wsChildren.Range("F13:H13").Value = 0.8
$wsChildren.Range("F14:H14").Value = 0.85
$wsChildren.Range("F15:H15").Value = 0.8
$wsChildren.Range("F16:H16").Value = 0.75
$wsChildren.Range("D18").Value = 0.19
$wsChildren.Range("D20").Value = 0.19
$wsChildren.Range("D22").Value = 0.19
$wsChildren.Range("D42:H42").Value = 0.5
$wsChildren.Range("D43:H43").Value = 0.63
$wsChildren.Range("D45:H45").Value = 0.8
$wsChildren.Range("D47:H47").Value = 0.76
$wsChildren.Range("E48").Value = 0
$wsChildren.Range("D49").Value = 0.88
$wsChildren.Range("E49").Value = 0
$wsChildren.Range("E2").Select()

# --- Programs for PW: effectiveness value updated ---
$wsPW = $wb.Worksheets.Item("Programs for PW")
$wsPW.Activate()
$wsPW.Range("D7:G7").Value = 0.59
$wsPW.Range("D7:G7").Select()

# --- Time trends was active before; Programs target population is now ---
$wsTarget = $wb.Worksheets.Item("Programs target population")
$wsTarget.Activate()
$wsTarget.Range("D9").Select()
